# DatosCreditoEmpresarialAmpliacion.xlsx - update test data for Ampliacion/Refinanciacion classes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 2 holds the sample data used to drive the form; refresh the client,
# pagare and monto figures used by the new test case.
$ws.Range("A2").Value = "23770571"          # Cliente
$ws.Range("H2").Value = "080-01-0456093"    # Pagare
$ws.Range("I2").Value = "2000"              # Monto

# Numero Propuesta / Resultado sample values are no longer populated.
$ws.Range("U2").ClearContents()
$ws.Range("V2").ClearContents()

# Resultado column header now reflects the ADN validation result.
$ws.Range("V1").Value = "Resultado con ADN"
